# correct the error in ppt
#
# 1) Fix the cached "datetimeFigureOut" date field shown on the Slide
#    Master and on every Slide Layout: 2020/3/16 -> 2020/4/19
# 2) Fix two typos on slide 5 ("上九" -> "上六" and "九五" -> "六五")

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $pf = $shp.PlaceholderFormat
            if ($pf.Type -eq 16) {
                # ppPlaceholderDate
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$p = $ppt.ActivePresentation

# --- 1) Slide Master + all Slide Layouts: refresh the stale date ---
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

Set-DatePlaceholderText $master.Shapes "2020/4/19"

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "2020/4/19"
}

# --- 2) Slide 5: correct two mis-typed hexagram line labels ---
$slide5 = $p.Slides.Item(5)
$shp = $slide5.Shapes.Item(5)
$tr = $shp.TextFrame.TextRange

$para1 = $tr.Paragraphs(1, 1)
$tr.Characters($para1.Start, 2).Text = "上六"

$para2 = $tr.Paragraphs(2, 1)
$tr.Characters($para2.Start, 2).Text = "六五"
